# "Testing beginners selection v3"
# Adds six new "Beginner" participants (rows 22-27) to Sheet1, wires up
# their e-mail hyperlinks, and updates the saved selection / window view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new participant data -------------------------------------------------
# (row, nummer, naam, email, ballroom/latin niveau, ballroom/latin lead-follow)
$people = @(
    @{ Row = 27; Nummer = 26; Naam = "Bibi";   Email = "marloes@gmail.com"; Niveau = "Beginner"; LeadFollow = "Follow" },
    @{ Row = 26; Nummer = 25; Naam = "Berit";  Email = "karin@gmail.com";   Niveau = "Beginner"; LeadFollow = "Follow" },
    @{ Row = 25; Nummer = 24; Naam = "Milan";  Email = "jasper@gmail.com"; Niveau = "Beginner"; LeadFollow = "Lead"   },
    @{ Row = 24; Nummer = 23; Naam = "Esther"; Email = "ilona@gmail.com";  Niveau = "Beginner"; LeadFollow = "Follow" },
    @{ Row = 23; Nummer = 22; Naam = "Bob";    Email = "harm@gmail.com";   Niveau = "Beginner"; LeadFollow = "Lead"   },
    @{ Row = 22; Nummer = 21; Naam = "Bjorn";  Email = "gydeon@gmail.com"; Niveau = "Beginner"; LeadFollow = "Lead"   }
)

# Fill the sheet bottom-up (27 -> 22), e-mail column before the name column,
# matching how the rows were originally typed in.
foreach ($p in $people) {
    $r = $p.Row
    $ws.Range("C$r").Value = $p.Email
    $ws.Range("B$r").Value = $p.Naam
    $ws.Range("A$r").Value = $p.Nummer
    $ws.Range("D$r").Value = $p.Niveau
    $ws.Range("E$r").Value = $p.Niveau
    $ws.Range("H$r").Value = $p.LeadFollow
    $ws.Range("I$r").Value = $p.LeadFollow
}

# --- hyperlinks -------------------------------------------------------------
# Added afterwards, in this particular order.
$hyperlinkRows = @(24, 25, 26, 27, 23, 22)
foreach ($r in $hyperlinkRows) {
    $email = ($people | Where-Object { $_.Row -eq $r }).Email
    $ws.Hyperlinks.Add($ws.Range("C$r"), "mailto:$email")
    $ws.Range("C$r").Style = "Hyperlink"
}

# --- selection / view --------------------------------------------------------
$ws.Range("G26").Select()

$win = $wb.Windows.Item(1)
$win.Top = 9600

Write-Output "done"
